# Replace the hyphen range separator with " to " in the "Expected" (column E)
# deaths-range values of the table, per the commit:
#   "Adapted table 2, new proposition for figure 2"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value  = "55,700 to 61,100"
$ws.Range("E3").Value  = "69,300 to 77,000"
$ws.Range("E4").Value  = "47,600 to 53,100"
$ws.Range("E5").Value  = "73,500 to 84,000"
$ws.Range("E6").Value  = "429,700 to 480,000"
$ws.Range("E7").Value  = "49,300 to 54,700"
$ws.Range("E8").Value  = "67,800 to 73,000"
$ws.Range("E9").Value  = "264,700 to 296,200"
$ws.Range("E10").Value = "64,700 to 70,900"
$ws.Range("E11").Value = "86,900 to 93,700"
$ws.Range("E12").Value = "400,500 to 441,000"
$ws.Range("E13").Value = "33,700 to 37,700"
$ws.Range("E15").Value = "210,300 to 238,400"
